# Update the exposure-sites table with the latest set of locations/sites.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @('Brighton', 'Brighton Beach  Brighton, VIC 3186', '26/12/20 12:00pm-1:00pm', 'Case attended beach', 'new'),
    @('Burwood East', 'Coles Burwood, Burwood Highway & Blackburn Road', '28/12/20 6.30pm - 7pm', 'Case shopped in store', 'new'),
    @('Burwood East', 'Kmart Burwood, 172 Burwood Highway', '28/12/20 6.15pm - 6.30pm', 'Case shopped in store', 'new'),
    @('Camberwell', 'Coles Middle Camberwell, 751 Riversdale Road', '28/12/20 12pm - 12.30pm', 'Case shopped in store', 'new'),
    @('Cheltenham', 'Cotton On, Southland Shopping Centre 1239 Nepean Hwy', '22/12/20 12.15pm - 12.45pm', 'Case visited venue', 'new'),
    @('Cheltenham', 'Myer, Southland Shopping Centre 1239 Nepean Hwy', '22/12/20 10.30am - 11am', 'Case visited venue', 'new'),
    @('Clayton', 'Kmart - M-City, 2107 Dandenong Road, Clayton', '30/12/20 7pm - 7.30pm', 'Case shopped at store', 'new'),
    @('Clayton', 'Woolworths - M-City, 2107 Dandenong Road, Clayton', '30/12/20 7.30pm - 745pm', 'Case shopped at store', 'new'),
    @('Dandenong', 'Kmart - Clayton Plaza, 2107 Dandenong Rd', '30/12/20 7pm - 7.30pm', 'Case shopped at store', 'old'),
    @('Dandenong', 'Woolworths - Clayton Plaza, 2107 Dandenong Rd', '30/12/20 7.30pm - 745pm', 'Case shopped at store', 'old'),
    @('Melbourne', 'Federation Square  Swanston & Flinders Streets, Melbourne VIC 3000', '23/12/2- 11:00pm-11:30pm', 'Case attended Federation Square', 'new'),
    @('Mentone', 'Bunnings Mentone  23-27 Nepean Hwy, Mentone VIC 3194', '23/12/20 11:00am-11:30am', 'Case shopped in store', 'new')
)

for ($i = 0; $i -lt $data.Length; $i++) {
    $row = $i + 2
    $rowData = $data[$i]
    for ($j = 0; $j -lt $rowData.Length; $j++) {
        $ws.Cells.Item($row, $j + 1).Value = $rowData[$j]
    }
}

$ws.Columns.Item(1).ColumnWidth = 10.833333333333334
$ws.Columns.Item(2).ColumnWidth = 54.833333333333336
$ws.Columns.Item(3).ColumnWidth = 23.5
$ws.Columns.Item(4).ColumnWidth = 26.166666666666668
